$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new test case row: AlertsNotificationManagement_TestClass / Yes
$ws.Range("A4").Value = "Core.AlertsandNotifications.AlertsNotificationManagement_TestClass"
$ws.Range("B4").Value = "Yes"

# Update the selection shown in the sheet view
$ws.Range("B6").Select()
